$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("QuantitativeMetrics")

# Update BLEU score (row 11)
$ws.Range("B11").Value = 0.1270730331643416

# Update Code BLEU score (row 12) and its descriptive note
$ws.Range("B12").Value = 0.3326171995683959
$ws.Range("C12").Value = "{'codebleu': 0.3326171995683959, 'ngram_match_score': 0.12707303316434163, 'weighted_ngram_match_score': 0.15302946474294168, 'syntax_match_score': 0.5384615384615384, 'dataflow_match_score': 0.5119047619047619}"

# Update Embeddings and Cosine similarity score (row 13)
$ws.Range("B13").Value = 0.8882556294488401
